$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this data block (rows 378-379),
# shifting the existing rows 378-397 down to 380-399.
$ws.Rows.Item(378).Resize(2).Insert()

# New row 378: Fukumoto, Primera, week of 2022-07-11
$ws.Range("A378").Value = 7
$ws.Range("B378").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C378").Value = "Ñuble"
$ws.Range("D378").Value = 44753
$ws.Range("E378").Value = 16
$ws.Range("F378").Value = "Fruta"
$ws.Range("G378").Value = 100102
$ws.Range("H378").Value = "Cítricos"
$ws.Range("I378").Value = 100102005
$ws.Range("J378").Value = "Naranja"
$ws.Range("K378").Value = "Fukumoto"
$ws.Range("L378").Value = "Primera"
$ws.Range("M378").Value = 160
$ws.Range("N378").Value = 7500
$ws.Range("O378").Value = 8000
$ws.Range("P378").Value = 7750
$ws.Range("Q378").Value = "$/bandeja 15 kilos granel"
$ws.Range("R378").Value = "Región de O'Higgins"
$ws.Range("S378").Value = 517
$ws.Range("T378").Value = 15

# New row 379: Fukumoto, Segunda, week of 2022-07-11
$ws.Range("A379").Value = 7
$ws.Range("B379").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C379").Value = "Ñuble"
$ws.Range("D379").Value = 44753
$ws.Range("E379").Value = 16
$ws.Range("F379").Value = "Fruta"
$ws.Range("G379").Value = 100102
$ws.Range("H379").Value = "Cítricos"
$ws.Range("I379").Value = 100102005
$ws.Range("J379").Value = "Naranja"
$ws.Range("K379").Value = "Fukumoto"
$ws.Range("L379").Value = "Segunda"
$ws.Range("M379").Value = 80
$ws.Range("N379").Value = 7000
$ws.Range("O379").Value = 7000
$ws.Range("P379").Value = 7000
$ws.Range("Q379").Value = "$/bandeja 15 kilos granel"
$ws.Range("R379").Value = "Región de O'Higgins"
$ws.Range("S379").Value = 467
$ws.Range("T379").Value = 15
